$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-0.03986657931951261, 1.166715395314321, 9.175165609364871, 3.029053583112202, 3.058341104199183, 52),
    @(0.3935572759351038, 0.9729011477362982, 6.967390270327681, 2.639581457414732, 2.636048747596285, 51),
    @(0.5845645237101601, 0.9612686620598623, 4.440372314603645, 2.107219095064309, 2.045067892746351, 50),
    @(0.4605577599719269, 1.029373326337189, 4.397813778901289, 2.097096511584836, 2.071314277662061, 41),
    @(0.5450649147966192, 1.227243952030447, 6.248613964862616, 2.499722777602072, 2.479899623702079, 31),
    @(0.4748548997082894, 1.211225564801436, 6.526816765525255, 2.554763543955733, 2.554677489448803, 29),
    @(0.411137621314042, 1.221296936513468, 6.896165483727119, 2.626055118181474, 2.643079106050209, 27),
    @(0.4034568762020432, 1.549620193073209, 9.580566432589528, 3.095249009787343, 3.15293505810293, 19),
    @(0.159026239978337, 1.962138362601939, 14.29507177400229, 3.780882406793723, 3.945508253329552, 12),
    @(-1.776008350391145, 2.037774147513916, 13.78833596834753, 3.713264866441327, 3.645910432883749, 5)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
